{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Title ---------------------------------------------------------\nparagraphs.items[0].insertText(\n  \"Science and Art: An Inseparable Union\",\n  Word.InsertLocation.replace\n);\n\n// --- 2. Author name -----------------------------------------------------\nparagraphs.items[1].insertText(\n  \"Amelia J. Masters\",\n  Word.InsertLocation.replace\n);\n\n// --- 3. Author email ------------------------------------------------------\nparagraphs.items[2].insertText(\n  \"amelia.masters@valid.edu\",\n  Word.InsertLocation.replace\n);\n\n// --- 4. Body paragraph (quantum -> science/art essay) --------------------\nconst bodyText =\n  \"As we venture through the world, our senses are met with an array of \" +\n  \"marvels that spark questions: How do we perceive objects in our \" +\n  \"environment? How do medications heal our bodies? Why does a tree shed \" +\n  \"its leaves in autumn? Exploring these queries leads us to a fascinating \" +\n  \"intersection where science and art harmoniously coexist. In this \" +\n  \"essay, we will traverse the extraordinary realm of the \" +\n  \"interrelatedness between science and art, revealing how these \" +\n  \"seemingly disparate disciplines converge to illuminate the profound \" +\n  \"mysteries of life\\u000b\\u000b\" +\n  \"In unraveling the intricate webs of scientific phenomena, researchers \" +\n  \"often draw inspiration from the beauty of art. The patterns and \" +\n  \"symmetries found in nature mirror the aesthetic principles guiding \" +\n  \"artistic expression. The spiral of a fern, the structure of a \" +\n  \"molecule, and the undulations of sound waves are just a few examples \" +\n  \"of the physical world's symphony that resonate with an artistic \" +\n  \"melody.\\u000b\\u000b\" +\n  \"Conversely, the methods of science serve as a tool for artists to \" +\n  \"convey their concepts visually and conceptually. The knowledge of \" +\n  \"color theory, perspective, and balance, which stems from scientific \" +\n  \"exploration and experimentation, empowers artists to craft \" +\n  \"masterpieces that explore the visible world's depth and complexity\";\n\nparagraphs.items[4].insertText(bodyText, Word.InsertLocation.replace);\n\n// --- 5. \"Summary\" heading is unchanged ------------------------------------\n\n// --- 6. Summary body paragraph --------------------------------------------\nconst summaryText =\n  \"Our exploration into the symbiotic relationship between science and \" +\n  \"art unveils how these domains, often perceived as distinct, are \" +\n  \"intertwined in a profound and enriching dance. Art, with its \" +\n  \"aesthetically pleasing forms and patterns, mirrors the structures and \" +\n  \"principles found in nature, while science provides tools and \" +\n  \"techniques that artists employ to present their ideas visually and \" +\n  \"conceptually. Through this harmonious convergence, we discover a \" +\n  \"universe where science illuminates the beauty of art, and art \" +\n  \"enhances our understanding of the natural world, enriching our \" +\n  \"perception and appreciation of the universe's wonders.\";\n\nparagraphs.items[6].insertText(summaryText, Word.InsertLocation.replace);\n\nawait context.sync();\n\n// --- 7. New empty paragraph appended at the very end of the body ---------\nbody.paragraphs.load(\"items\");\nawait context.sync();\nconst lastParaList = body.paragraphs;\nconst lastPara = lastParaList.items[lastParaList.items.length - 1];\nlastPara.insertParagraph(\"\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d (ActiveDocument) are pre-seeded by the harness.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParaText($Paragraph, $NewText) {\n    $r = $Paragraph.Range\n    # Exclude the trailing paragraph mark so the paragraph's formatting\n    # (pPr / the run's rPr) survives the rewrite.\n    $r.End = $r.End - 1\n    $r.Text = $NewText\n}\n\n# --- 1. Title --------------------------------------------------------------\nSet-ParaText $d.Paragraphs.Item(1) \"Science and Art: An Inseparable Union\"\n\n# --- 2. Author name ----------------------------------------------------------\nSet-ParaText $d.Paragraphs.Item(2) \"Amelia J. Masters\"\n\n# --- 3. Author email -----------------------------------------------------------\nSet-ParaText $d.Paragraphs.Item(3) \"amelia.masters@valid.edu\"\n\n# --- 4. Body paragraph (quantum -> science/art essay) -----------------------\n$bodyText = \"As we venture through the world, our senses are met with an array of \" +\n  \"marvels that spark questions: How do we perceive objects in our \" +\n  \"environment? How do medications heal our bodies? Why does a tree shed \" +\n  \"its leaves in autumn? Exploring these queries leads us to a fascinating \" +\n  \"intersection where science and art harmoniously coexist. In this \" +\n  \"essay, we will traverse the extraordinary realm of the \" +\n  \"interrelatedness between science and art, revealing how these \" +\n  \"seemingly disparate disciplines converge to illuminate the profound \" +\n  \"mysteries of life\" + [char]11 + [char]11 +\n  \"In unraveling the intricate webs of scientific phenomena, researchers \" +\n  \"often draw inspiration from the beauty of art. The patterns and \" +\n  \"symmetries found in nature mirror the aesthetic principles guiding \" +\n  \"artistic expression. The spiral of a fern, the structure of a \" +\n  \"molecule, and the undulations of sound waves are just a few examples \" +\n  \"of the physical world's symphony that resonate with an artistic \" +\n  \"melody.\" + [char]11 + [char]11 +\n  \"Conversely, the methods of science serve as a tool for artists to \" +\n  \"convey their concepts visually and conceptually. The knowledge of \" +\n  \"color theory, perspective, and balance, which stems from scientific \" +\n  \"exploration and experimentation, empowers artists to craft \" +\n  \"masterpieces that explore the visible world's depth and complexity\"\n\nSet-ParaText $d.Paragraphs.Item(5) $bodyText\n\n# --- 5. \"Summary\" heading is unchanged ---------------------------------------\n\n# --- 6. Summary body paragraph ----------------------------------------------\n$summaryText = \"Our exploration into the symbiotic relationship between science and \" +\n  \"art unveils how these domains, often perceived as distinct, are \" +\n  \"intertwined in a profound and enriching dance. Art, with its \" +\n  \"aesthetically pleasing forms and patterns, mirrors the structures and \" +\n  \"principles found in nature, while science provides tools and \" +\n  \"techniques that artists employ to present their ideas visually and \" +\n  \"conceptually. Through this harmonious convergence, we discover a \" +\n  \"universe where science illuminates the beauty of art, and art \" +\n  \"enhances our understanding of the natural world, enriching our \" +\n  \"perception and appreciation of the universe's wonders.\"\n\nSet-ParaText $d.Paragraphs.Item(7) $summaryText\n\n# --- 7. New empty paragraph appended at the very end of the body ------------\n$endRange = $d.Content\n$endRange.Collapse(0) | Out-Null\n$endRange.InsertParagraphAfter() | Out-Null\n"}
